$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text would otherwise be auto-parsed as a number
# are forced to remain text (matching the original inlineStr cell type)
# by temporarily applying a text number format, then restoring the default
# "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").Value = "57.939.96"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.128.18"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.419"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.74%  "
$ws.Range("D13").Value = "3.662.02"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "58.009.37"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.64%  "
$ws.Range("D18").Value = "3.126.22"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.94%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("E36").Value = "  +3.35%  "
$ws.Range("E37").Value = "  +7.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0676"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").Value = "2.548.67"
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.60%  "
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0988"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.750"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.49%  "
